$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first fixture row (Real Madrid CF - Osasuna / 27/10/2021);
# this shifts every subsequent row up by one.
$ws.Rows("1:1").Delete()

# The "Real Madrid CF - Rayo Vallecano" fixture (now row 2) moved from
# 07/11/2021 to 06/11/2021. Force the cell to stay text (matching the
# existing text-typed date cells) instead of being auto-parsed as a date,
# then restore the default (unstyled) cell formatting.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "06/11/2021"
$ws.Range("B2").Style = "Normal"

# Append the new fixture at the end of the table (now row 10).
$ws.Range("A10").Value = "Real Madrid CF - Deportivo Alavés"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "20/02/2022"
$ws.Range("B10").Style = "Normal"
